$d = $word.ActiveDocument

# Locate the paragraph that contains the "Ver no Jupiter ..." text. It is
# immediately preceded by a blank paragraph and immediately followed by the
# "(c) 2020 ... Creative Commons Attribution" paragraph. All three paragraphs
# (the blank one, the "Ver no Jupiter" one, and the copyright one) need to be
# removed, leaving the "LOQ4031: ..." paragraph followed directly by the
# (already-present) blank paragraph that precedes the page-break paragraph.

$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Ver no Jupiter*") {
        $targetIdx = $i
        break
    }
}

if ($targetIdx -gt 1) {
    $pBefore = $d.Paragraphs.Item($targetIdx - 1)
    $pAfter = $d.Paragraphs.Item($targetIdx + 1)

    $delStart = $pBefore.Range.Start
    $delEnd = $pAfter.Range.End

    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
